$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("A20").Value2 = "TestLabel"
$src = $ws.Range("A1")
$dst = $ws.Range("A20")
$dst.Font.Bold = $src.Font.Bold
$dst.HorizontalAlignment = $src.HorizontalAlignment
$dst.VerticalAlignment = $src.VerticalAlignment
$dst.Borders.LineStyle = $src.Borders.LineStyle
Write-Host "done"
